# app/celery/tasks/tag/files/tags.xlsx
# "add tv-8, add tasks get pressure, get min temperature"
#
# Append four new tag rows (title/description pairs) below the existing
# data, add a trailing blank styled row, widen the two columns to fit the
# new (longer) content, and leave the sheet selection on E3 - matching the
# author's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new tag rows -----------------------------------------------------
$ws.Range("A8").Value = "AverageAbsolutePressurePerDay"
$ws.Range("B8").Value = "Среднее абсолютное давление за сутки"

$ws.Range("A9").Value = "TemperatureMinPerDay"
$ws.Range("B9").Value = "Минимальная температура за сутки"

$ws.Range("A10").Value = "VolumeMaxPerDay"
$ws.Range("B10").Value = "Максимальный объем за сутки"

$ws.Range("A11").Value = "VolumeForwardFixDay"
$ws.Range("B11").Value = "Объем в прямом направлении на начало суток"

# Carry the existing row style (font "Times New Roman", cellXf index 1)
# onto the freshly written cells - same as every other row on the sheet -
# plus a new trailing blank row (A12:B12), by cloning formats from the
# last pre-existing row.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- column widths ------------------------------------------------------
# Widen the columns so the new, longer strings fit (matches the author's
# "best fit" widths of 28 / 41.81640625 characters - the engine quantizes
# the stored width to 1/6-character steps, so 41 is the closest input that
# rounds to the saved 41.83333... width nearest 41.81640625).
$ws.Columns.Item(1).ColumnWidth = 27.166666666666668
$ws.Columns.Item(2).ColumnWidth = 41

# --- selection ------------------------------------------------------
$ws.Range("E3").Select()
